$wb = $excel.ActiveWorkbook

# --- Sheet 1: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ25538881"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 5647.544512168372
$ws.Cells.Item(2, 3).Value = 0.03432310885369553
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -7574.809392407406
$ws.Cells.Item(3, 3).Value = 0.003415386611786744
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -5197.132716366207
$ws.Cells.Item(4, 3).Value = 0.03336673489035158
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -4942.66039963905
$ws.Cells.Item(5, 3).Value = 0.04287901750179909
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 731.451744136182
$ws.Cells.Item(6, 3).Value = 0.01827180881499101
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = 162.224306315782
$ws.Cells.Item(7, 3).Value = 0.5129828433246478
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 150.6274958157008
$ws.Cells.Item(8, 3).Value = 0.0505878142367687
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1198.386116051623
$ws.Cells.Item(9, 3).Value = 0.0000000003145741265608967
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -13.24571079802735
$ws.Cells.Item(10, 3).Value = 0.1100784490350273
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -192.9384582915875
$ws.Cells.Item(11, 3).Value = 0.01368743681284751
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 738.4571545755616
$ws.Cells.Item(12, 3).Value = 0.0000000000000000000000000000003068378122710816
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.01742097655690067
$ws.Cells.Item(13, 3).Value = 0.5919110479774294
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = 0.00005214441441001393
$ws.Cells.Item(14, 3).Value = 0.3978097168806588
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -4.409317097707444
$ws.Cells.Item(15, 3).Value = 0.5637522818156443
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 39.31293857636881
$ws.Cells.Item(16, 3).Value = 0.000001205254940222836
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -60.59605326519105
$ws.Cells.Item(17, 3).Value = 0.9562214962920521
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 460.2071264478043
$ws.Cells.Item(18, 3).Value = 0.5651432000391414
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 5274.158472510791
$ws.Cells.Item(19, 3).Value = 0.2854388144880188

# --- Sheet 2: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ25803405"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 7746.004600637669
$ws.Cells.Item(2, 3).Value = 0.003806113442815872
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -5183.004052760492
$ws.Cells.Item(3, 3).Value = 0.04742987675525355
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -5069.023839372901
$ws.Cells.Item(4, 3).Value = 0.03865160275164621
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -4859.226204019713
$ws.Cells.Item(5, 3).Value = 0.04729510209658014
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 670.7882760904707
$ws.Cells.Item(6, 3).Value = 0.03082663223014452
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -211.0817643205829
$ws.Cells.Item(7, 3).Value = 0.3919386267963497
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 57.07851913500264
$ws.Cells.Item(8, 3).Value = 0.4658879646377384
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1202.707996162474
$ws.Cells.Item(9, 3).Value = 0.0000000002811698601214166
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -23.04650721862184
$ws.Cells.Item(10, 3).Value = 0.005658423374751379
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -183.7494578060234
$ws.Cells.Item(11, 3).Value = 0.01813654938444839
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 693.6929108005293
$ws.Cells.Item(12, 3).Value = 0.000000000000000000000000000323472439001115
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.01325454033477088
$ws.Cells.Item(13, 3).Value = 0.6912741126742425
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = -0.0000231883319787489
$ws.Cells.Item(14, 3).Value = 0.7114987665844192
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -1.430816555076977
$ws.Cells.Item(15, 3).Value = 0.8515533795506866
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 32.05976206777068
$ws.Cells.Item(16, 3).Value = 0.0000543832537730159
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -607.003012011224
$ws.Cells.Item(17, 3).Value = 0.5816150210705746
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = -325.6646379531379
$ws.Cells.Item(18, 3).Value = 0.6831737942787628
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 1226.288341961202
$ws.Cells.Item(19, 3).Value = 0.7993411936340062

# --- Sheet 3: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ26054529"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 2724.53909888409
$ws.Cells.Item(2, 3).Value = 0.4553739081838598
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -1250.217367766923
$ws.Cells.Item(3, 3).Value = 0.7253937429970994
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 319.0873732480509
$ws.Cells.Item(4, 3).Value = 0.9265994213027539
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 595.0904895723206
$ws.Cells.Item(5, 3).Value = 0.8635331351878983
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 915.6248568993026
$ws.Cells.Item(6, 3).Value = 0.003280811552561296
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -14.16107496576836
$ws.Cells.Item(7, 3).Value = 0.9542119201333457
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 65.7937493428467
$ws.Cells.Item(8, 3).Value = 0.4006517501565446
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1205.498001781406
$ws.Cells.Item(9, 3).Value = 0.0000000002745757261276916
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -12.66492896542794
$ws.Cells.Item(10, 3).Value = 0.133080753125965
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -150.1269626967944
$ws.Cells.Item(11, 3).Value = 0.05365141068612737
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 620.2063766865655
$ws.Cells.Item(12, 3).Value = 0.00000000000000000000009065582022740617
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.01058218809172087
$ws.Cells.Item(13, 3).Value = 0.7436199754389473
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = 0.00002208855661645376
$ws.Cells.Item(14, 3).Value = 0.7178551833104192
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -13.1268238037241
$ws.Cells.Item(15, 3).Value = 0.08928814278480698
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 26.54324899128483
$ws.Cells.Item(16, 3).Value = 0.001065872100891412
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -272.76894501259
$ws.Cells.Item(17, 3).Value = 0.8045779659394061
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 104.0023192433221
$ws.Cells.Item(18, 3).Value = 0.8997289522610444
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 5607.808861356883
$ws.Cells.Item(19, 3).Value = 0.2485638194026963

# --- Sheet 4: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ26313895"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 6717.529522141865
$ws.Cells.Item(2, 3).Value = 0.01277558940785771
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -6369.073687684883
$ws.Cells.Item(3, 3).Value = 0.01395987472831499
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -4845.474454099205
$ws.Cells.Item(4, 3).Value = 0.04813183096804059
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -4696.351464471632
$ws.Cells.Item(5, 3).Value = 0.05528151948842362
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 616.6788615438555
$ws.Cells.Item(6, 3).Value = 0.04898983159998931
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -290.5790261607431
$ws.Cells.Item(7, 3).Value = 0.2485416932637795
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 24.4456060955863
$ws.Cells.Item(8, 3).Value = 0.7531348530748443
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1212.214572793829
$ws.Cells.Item(9, 3).Value = 0.0000000002466460480292611
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -17.27424407184199
$ws.Cells.Item(10, 3).Value = 0.0370643144675265
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -119.448254600765
$ws.Cells.Item(11, 3).Value = 0.1264278537949756
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 645.2026828694927
$ws.Cells.Item(12, 3).Value = 0.000000000000000000000003077259579513356
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.01668732562483845
$ws.Cells.Item(13, 3).Value = 0.6014559657607096
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = -0.000001103456406903955
$ws.Cells.Item(14, 3).Value = 0.9854676051286689
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -4.833185045363391
$ws.Cells.Item(15, 3).Value = 0.525390915671579
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 40.66106513955712
$ws.Cells.Item(16, 3).Value = 0.000003586117078274811
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -373.7933360945322
$ws.Cells.Item(17, 3).Value = 0.7346934356994528
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 335.5636266971835
$ws.Cells.Item(18, 3).Value = 0.6796774011856754
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 2490.356387904053
$ws.Cells.Item(19, 3).Value = 0.6079129875208034

# --- Sheet 5: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ26585732"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 7547.082890011134
$ws.Cells.Item(2, 3).Value = 0.004628346822745491
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -6640.740476463921
$ws.Cells.Item(3, 3).Value = 0.009824745322862571
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -5118.103247854749
$ws.Cells.Item(4, 3).Value = 0.03621071563936388
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -4920.394177090681
$ws.Cells.Item(5, 3).Value = 0.04388270334864274
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 655.8506677030196
$ws.Cells.Item(6, 3).Value = 0.03366981675849895
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -236.9138491273626
$ws.Cells.Item(7, 3).Value = 0.3402121641298428
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 112.0937511879716
$ws.Cells.Item(8, 3).Value = 0.1463401363423115
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1323.500957577215
$ws.Cells.Item(9, 3).Value = 0.000000000004191669654022238
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -14.10786136643806
$ws.Cells.Item(10, 3).Value = 0.08812506218422853
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -196.3456201026291
$ws.Cells.Item(11, 3).Value = 0.01187389193420413
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 629.3981717768165
$ws.Cells.Item(12, 3).Value = 0.00000000000000000000001621170529109463
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.0127942663050716
$ws.Cells.Item(13, 3).Value = 0.6887492583478361
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = 0.00004591797163832441
$ws.Cells.Item(14, 3).Value = 0.4523419974618065
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -9.863149462773702
$ws.Cells.Item(15, 3).Value = 0.1928077172643367
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 35.89136068615463
$ws.Cells.Item(16, 3).Value = 0.000007821032718143633
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -1308.128821845608
$ws.Cells.Item(17, 3).Value = 0.2356442909574183
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 100.5245285396899
$ws.Cells.Item(18, 3).Value = 0.9018470758260622
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 4687.111899113117
$ws.Cells.Item(19, 3).Value = 0.3368215245743306

# --- Sheet 6: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ26824756"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 14121.01136454242
$ws.Cells.Item(2, 3).Value = 0.0000859241182350305
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -11739.70622758823
$ws.Cells.Item(3, 3).Value = 0.001052761829722274
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -10709.49569759975
$ws.Cells.Item(4, 3).Value = 0.001901462744623056
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -10317.60222277377
$ws.Cells.Item(5, 3).Value = 0.002766853009426637
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 589.8596001810047
$ws.Cells.Item(6, 3).Value = 0.05652363579608755
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -254.5932408895038
$ws.Cells.Item(7, 3).Value = 0.3056364893486381
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 103.7690512581148
$ws.Cells.Item(8, 3).Value = 0.1771328275222099
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1042.501698797422
$ws.Cells.Item(9, 3).Value = 0.00000003767717003966091
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -13.51388343588172
$ws.Cells.Item(10, 3).Value = 0.0981637109204173
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -267.9923071027872
$ws.Cells.Item(11, 3).Value = 0.0006826192777913909
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 675.7050233639511
$ws.Cells.Item(12, 3).Value = 0.000000000000000000000000004570904501723446
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.009662326169680966
$ws.Cells.Item(13, 3).Value = 0.7598617445940882
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = 0.000008205143126691895
$ws.Cells.Item(14, 3).Value = 0.8925492721179193
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -10.2594659650664
$ws.Cells.Item(15, 3).Value = 0.1724267586307751
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 24.94298803521974
$ws.Cells.Item(16, 3).Value = 0.002760346810659689
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -1228.885565551523
$ws.Cells.Item(17, 3).Value = 0.2627526354995405
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = -63.16081388924911
$ws.Cells.Item(18, 3).Value = 0.9376176654076008
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = -812.5910942862133
$ws.Cells.Item(19, 3).Value = 0.8651946622295621

# --- Sheet 7: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ27062304"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 7370.855798006472
$ws.Cells.Item(2, 3).Value = 0.005283598422928783
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -7025.461266294036
$ws.Cells.Item(3, 3).Value = 0.005981942546836079
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -5188.880624612727
$ws.Cells.Item(4, 3).Value = 0.03211087712469719
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -4929.330297590566
$ws.Cells.Item(5, 3).Value = 0.04162815803699558
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 631.5403941444912
$ws.Cells.Item(6, 3).Value = 0.04035525995629005
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -10.13703821677956
$ws.Cells.Item(7, 3).Value = 0.9666476788956049
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 140.8873679183746
$ws.Cells.Item(8, 3).Value = 0.06351650090953334
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -880.9936945537156
$ws.Cells.Item(9, 3).Value = 0.000003312506368602994
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -10.56235800391658
$ws.Cells.Item(10, 3).Value = 0.1963931832466881
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -115.5891193257134
$ws.Cells.Item(11, 3).Value = 0.134869241621737
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 629.1355027261493
$ws.Cells.Item(12, 3).Value = 0.000000000000000000000001089601550555827
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.007728724162998019
$ws.Cells.Item(13, 3).Value = 0.8069648544378809
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = 0.0000223072119710135
$ws.Cells.Item(14, 3).Value = 0.7099211913615275
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -8.715495407117222
$ws.Cells.Item(15, 3).Value = 0.2425031714317525
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 29.05252969681494
$ws.Cells.Item(16, 3).Value = 0.0002946388382551976
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -707.2997241564924
$ws.Cells.Item(17, 3).Value = 0.5152934230616022
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = -210.8640047158188
$ws.Cells.Item(18, 3).Value = 0.793633381003481
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 1788.061095113567
$ws.Cells.Item(19, 3).Value = 0.7140499639168795

# --- Sheet 8: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ27328602"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 7405.537600870319
$ws.Cells.Item(2, 3).Value = 0.004900317877676889
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -7112.519665815447
$ws.Cells.Item(3, 3).Value = 0.005191407631390582
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -5202.884980402171
$ws.Cells.Item(4, 3).Value = 0.03044898010808985
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -4998.980529545478
$ws.Cells.Item(5, 3).Value = 0.03744045427014081
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 860.731107022893
$ws.Cells.Item(6, 3).Value = 0.00494849329149876
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -8.373795826970024
$ws.Cells.Item(7, 3).Value = 0.9724087298333515
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 234.206379107817
$ws.Cells.Item(8, 3).Value = 0.002438037548249375
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1020.016453486038
$ws.Cells.Item(9, 3).Value = 0.00000005530055665504536
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -17.21735037476993
$ws.Cells.Item(10, 3).Value = 0.03333081060019757
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -187.6072267281807
$ws.Cells.Item(11, 3).Value = 0.01447068388257372
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 650.5459190632464
$ws.Cells.Item(12, 3).Value = 0.0000000000000000000000001043429159087501
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.05137944978188086
$ws.Cells.Item(13, 3).Value = 0.1034993057016179
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = 0.00005984423092021661
$ws.Cells.Item(14, 3).Value = 0.3184188384008496
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -10.63581188562917
$ws.Cells.Item(15, 3).Value = 0.1532634714170395
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 30.57247776216976
$ws.Cells.Item(16, 3).Value = 0.000217443342357554
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -265.7218259677652
$ws.Cells.Item(17, 3).Value = 0.8068596855735103
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 646.7875154519718
$ws.Cells.Item(18, 3).Value = 0.4123454518533695
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = -894.5730588550368
$ws.Cells.Item(19, 3).Value = 0.8509818308844542

# --- Sheet 9: rerun dist commute with harmonised education ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ27572637"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 8116.304817696619
$ws.Cells.Item(2, 3).Value = 0.002338777500836103
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = -6405.51486308447
$ws.Cells.Item(3, 3).Value = 0.01299544011810717
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -5072.836687512365
$ws.Cells.Item(4, 3).Value = 0.03786419446378841
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -4853.119700944904
$ws.Cells.Item(5, 3).Value = 0.0468348043746788
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = 930.6129616198295
$ws.Cells.Item(6, 3).Value = 0.002871848970378708
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -78.97676916406373
$ws.Cells.Item(7, 3).Value = 0.7496353050187609
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 66.28469524007747
$ws.Cells.Item(8, 3).Value = 0.3935684242641793
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -1074.61279970195
$ws.Cells.Item(9, 3).Value = 0.00000001692588777585523
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -20.70595929417399
$ws.Cells.Item(10, 3).Value = 0.01215766347892219
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = -129.0814409358699
$ws.Cells.Item(11, 3).Value = 0.09379328476217241
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 612.3916061335287
$ws.Cells.Item(12, 3).Value = 0.00000000000000000000008167933519086865
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = 0.01603145871785329
$ws.Cells.Item(13, 3).Value = 0.6192860412174002
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = -0.00003808707446817847
$ws.Cells.Item(14, 3).Value = 0.5361923562879427
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -8.654080409955926
$ws.Cells.Item(15, 3).Value = 0.2515195878695601
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 29.8331442744315
$ws.Cells.Item(16, 3).Value = 0.0001803226359936902
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -773.6811048726299
$ws.Cells.Item(17, 3).Value = 0.4789358456931013
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = -275.890081840791
$ws.Cells.Item(18, 3).Value = 0.7330963182553552
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 8348.61670336093
$ws.Cells.Item(19, 3).Value = 0.09049827250104589

Write-Output "edit complete"
